$d = $word.ActiveDocument

# The state-chart drawing is the document's single inline lockedCanvas
# drawing. We read its OOXML, make the three edits from the diff
# (resize the drawing, renumber the docPr id, and add the new dashed
# connector shape before the closing lockedCanvas tag), then push the
# edited XML back over the same range.

$shp = $d.InlineShapes.Item(1)
$r = $shp.Range
$xml = $r.WordOpenXML

# Exporting a Range to WordOpenXML stamps a fresh w14:paraId/textId onto
# the paragraph so it can stand alone as a well-formed fragment; strip
# that back out since the real paragraph never had one and we don't want
# to invent one just because we had to round-trip through OOXML text.
$xml = $xml.Replace(' w14:paraId="00000001" w14:textId="77777777"', '')

# 1) Resize the drawing extent.
$xml = $xml.Replace('<wp:extent cx="4933950" cy="3152775"/>', '<wp:extent cx="4429125" cy="3009900"/>')

# 2) Renumber the drawing's docPr id.
$xml = $xml.Replace('<wp:docPr id="2" name="Object 1"/>', '<wp:docPr id="1" name="Object 1"/>')

# 3) Add the new dashed "Straight Connector 2" shape as the last shape
#    in the locked canvas, right before it closes.
$newConnector = '<a:cxnSp><a:nvCxnSpPr><a:cNvPr id="3" name="Straight Connector 2"/><a:cNvCxnSpPr/></a:nvCxnSpPr><a:spPr><a:xfrm><a:off x="5867400" y="533400"/><a:ext cx="76200" cy="6172200"/></a:xfrm><a:prstGeom prst="line"><a:avLst/></a:prstGeom><a:ln><a:prstDash val="dash"/></a:ln></a:spPr><a:style><a:lnRef idx="2"><a:schemeClr val="accent1"/></a:lnRef><a:fillRef idx="0"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="1"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="tx1"/></a:fontRef></a:style></a:cxnSp>'
$xml = $xml.Replace('</lc:lockedCanvas>', $newConnector + '</lc:lockedCanvas>')

$ins = $r.Duplicate
$ins.Collapse(1)
[void]$ins.InsertXML($xml)
